$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting old rows 78-79 down to 79-80
$ws.Rows.Item(78).Insert()

# Fill in the new row 78 with data
$ws.Cells.Item(78, 1).Value = 11
$ws.Cells.Item(78, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(78, 3).Value = "Bíobío"
$ws.Cells.Item(78, 4).Value = 44516
$ws.Cells.Item(78, 5).Value = 8
$ws.Cells.Item(78, 6).Value = "Fruta"
$ws.Cells.Item(78, 7).Value = 100108
$ws.Cells.Item(78, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(78, 9).Value = 100108002
$ws.Cells.Item(78, 10).Value = "Mango"
$ws.Cells.Item(78, 11).Value = "Sin especificar"
$ws.Cells.Item(78, 12).Value = "Primera"
$ws.Cells.Item(78, 13).Value = 400
$ws.Cells.Item(78, 14).Value = 7500
$ws.Cells.Item(78, 15).Value = 8000
$ws.Cells.Item(78, 16).Value = 7750
$ws.Cells.Item(78, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(78, 18).Value = "Perú"
$ws.Cells.Item(78, 19).Value = 1938
$ws.Cells.Item(78, 20).Value = 4
